$d = $word.ActiveDocument

# 1. Replace first sentence text
$d.Content.Find.Execute(
    "Hopefully I can do a similar operation for the putImageData method call",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "The same idea was applied to the putImageData call in the view", 2)

# 2. Replace second sentence text (the one that keeps the bookmark)
$d.Content.Find.Execute(
    "I will also look into what processes are the most expensive when adjusting an image" + [char]8217 + "s colour",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "The most expensive process is now the garbage collector, followed by definitions of Color objects: this leads me to believe that the best way to improve the performance of my program now is to reduce the number of objects declared and/or make their instantiation more efficient", 2)

# 3. Find the paragraph that now contains the replaced text (with bookmark) and
#    add a brand new list paragraph after it with the original sentence text.
$target = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "The most expensive process is now the garbage collector*") {
        $target = $p
        break
    }
}

$newPara = $d.Paragraphs.Add($target.Range)
$newPara.Range.Text = "I will also look into what processes are the most expensive when adjusting an image" + [char]8217 + "s colour"
